# Daily update at 8 AM UTC
# Adds the next day's row (45662 -> 2025-01-05) to the "Wins Over Time" log
# and moves the "last row" date formatting down from the old last row (74)
# to the new last row (75).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 74 was previously the last row and carried the distinct "last row"
# date format (YYYY-MM-DD, no time). Since it's no longer the last row,
# give it the regular date format used by all the other data rows.
$ws.Range("A74").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 75.
$ws.Range("A75").Value = 45662
$ws.Range("B75").Value = 173
$ws.Range("C75").Value = 174
$ws.Range("D75").Value = 175

# The new last row takes on the distinct "last row" date format.
$ws.Range("A75").NumberFormat = "YYYY-MM-DD"
